# Auto Update Data
# - Refresh the "更新日期" (last-updated) timestamp stamped in A1.
# - The task list rows were re-ordered/updated upstream: for rows 310-360
#   each row now shows the operator/stage/requirement text that used to sit
#   in the row directly below it (a one-row "shift up"). Row 361 (the last
#   entry, 行箸 / 3-2) is unchanged, so it ends up duplicated onto row 360
#   as well as remaining on row 361 itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the A/B/C (干员/关卡/内容) columns for rows 310..360 up by one row,
# pulling each row's new content from the row below it. Walking top-down is
# safe here because we always read row ($r + 1) before it is ever written.
for ($r = 310; $r -le 360; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($r + 1, 1).Value2
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r + 1, 2).Value2
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item($r + 1, 3).Value2
}

# Bump the "last updated" banner in A1 to match the new export timestamp.
$ws.Range("A1").Value2 = "更新日期：2025.01.31 13:17:43"
